$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20 - this shifts existing rows 20..25 down to 21..26
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with the latest week's data (same template as the
# surrounding Chirimoya rows, new date + values per the diff)
$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value = "Ñuble"
$ws.Cells.Item(20, 4).Value = 45205
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100107
$ws.Cells.Item(20, 8).Value = "Otros"
$ws.Cells.Item(20, 9).Value = 100107002
$ws.Cells.Item(20, 10).Value = "Chirimoya"
$ws.Cells.Item(20, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 30
$ws.Cells.Item(20, 14).Value = 22000
$ws.Cells.Item(20, 15).Value = 22000
$ws.Cells.Item(20, 16).Value = 22000
$ws.Cells.Item(20, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 19).Value = 2200
$ws.Cells.Item(20, 20).Value = 10

$wb.Save()
